$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "87.837.70"
$ws.Range("E2").Value = "  -5.62%  "

$ws.Range("D3").Value = "3.089.44"
$ws.Range("E3").Value = "  -7.08%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.81"
$ws.Range("E5").Value = "  -2.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "631.39"
$ws.Range("E6").Value = "  +0.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.388"
$ws.Range("E7").Value = "  -6.99%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.721"
$ws.Range("E8").Value = "  +1.48%  "

$ws.Range("E9").Value = "  +0.28%  "

$ws.Range("D10").Value = "3.085.78"
$ws.Range("E10").Value = "  -7.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.552"
$ws.Range("E11").Value = "  -7.32%  "

$ws.Range("E12").Value = "  -1.81%  "

$ws.Range("E13").Value = "  -8.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.25"
$ws.Range("E14").Value = "  -2.87%  "

$ws.Range("D15").Value = "87.667.91"
$ws.Range("E15").Value = "  -5.08%  "

$ws.Range("D16").Value = "3.662.43"
$ws.Range("E16").Value = "  -6.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "31.97"
$ws.Range("E17").Value = "  -8.41%  "

$ws.Range("D18").Value = "3.093.47"
$ws.Range("E18").Value = "  -6.73%  "

$ws.Range("E19").Value = "  -1.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000211"
$ws.Range("E20").Value = "  +7.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.04"
$ws.Range("E21").Value = "  -8.79%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "422.05"
$ws.Range("E22").Value = "  -4.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.29"
$ws.Range("E23").Value = "  -8.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.85"
$ws.Range("E24").Value = "  -9.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.24"
$ws.Range("E25").Value = "  -3.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.54"
$ws.Range("E26").Value = "  -8.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "78.87"
$ws.Range("E27").Value = "  +2.86%  "

$ws.Range("D28").Value = "3.250.80"
$ws.Range("E28").Value = "  -7.68%  "

$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("B30").Value = "Cronos"
$ws.Range("C30").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.156"
$ws.Range("E30").Value = "  -13.92%  "

$ws.Range("B31").Value = "dogwifhat"
$ws.Range("C31").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.95"
$ws.Range("E31").Value = "  +4.27%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.11"
$ws.Range("E32").Value = "  -8.70%  "

$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "505.28"
$ws.Range("E33").Value = "  -11.28%  "

$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.74"
$ws.Range("E34").Value = "  -9.09%  "

$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.26"
$ws.Range("E35").Value = "  -6.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.81"
$ws.Range("E36").Value = "  -6.38%  "

$ws.Range("B37").Value = "WhiteBITCoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "22.20"
$ws.Range("E37").Value = "  -0.99%  "

$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "21.59"
$ws.Range("E38").Value = "  -5.90%  "

$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.41%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.123"
$ws.Range("E40").Value = "  -7.53%  "

$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.84"
$ws.Range("E42").Value = "  -8.90%  "

$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.359"
$ws.Range("E43").Value = "  -10.53%  "

$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "143.96"
$ws.Range("E44").Value = "  -4.62%  "

$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.26"
$ws.Range("E45").Value = "  -1.35%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.123"
$ws.Range("E46").Value = "  -5.90%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "163.95"
$ws.Range("E47").Value = "  -10.66%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.712"
$ws.Range("E48").Value = "  -3.76%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.09"
$ws.Range("E49").Value = "  -5.66%  "

$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.588"
$ws.Range("E50").Value = "  -7.75%  "

$ws.Range("B51").Value = "Binance-PegBSC-USD"
$ws.Range("C51").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.545"
$ws.Range("E51").Value = "  -45.63%  "
